$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" record (originally row 26) entirely - rows below shift up.
$ws.Rows.Item(26).Delete()

# Remove the "SC 92" record (originally row 28, now row 27 after the prior delete).
$ws.Rows.Item(27).Delete()

# --- Update column F values to reflect the new "missing data" pattern ---

# RM 8 (row 3): previously missing, now has a value.
$ws.Range("F3").Value = 17.64

# RM 14 (row 5): previously had a value, now missing.
$ws.Range("F5").Value = ""

# RM 135 (row 21): previously missing, now has a value.
$ws.Range("F21").Value = 16.58

# RM 140 (row 23): previously had a value, now missing.
$ws.Range("F23").Value = ""

# SC 5 (row 26, formerly row 27 before the row deletions): previously a single
# space placeholder, now has a numeric value.
$ws.Range("F26").Value = 17.38

# SC 193 (row 32, formerly row 34 before the row deletions): previously missing,
# now has a value.
$ws.Range("F32").Value = 17.39

$wb.Save()
